# 9 testcase for E2E hydroflask
#
# Insert two new DataSet rows (SKU-CFX001 / SKU-T20CPB001) into the "E2E"
# sheet right after the existing row 13, pushing every row below (14-38)
# down by two, and refresh the gift-card redemption code used by the
# "Full_RedeemGiftcard" test case (now on row 38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("E2E")

# Insert two blank rows before the current row 14. Excel copies the
# formatting of the row immediately above (row 13), which already carries
# the styles (s="33" on column AB, s="5" on AC, s="2" on H/I) that the new
# rows need.
$ws.Range("A14:A15").EntireRow.Insert()

# New row 14: Wide Mouth Flex Sip (tm) Lid bundle test case
$ws.Range("A14").Value = "SKU-CFX001"
$ws.Range("AB14").Value = "Wide Mouth Flex Sip" + [char]0x2122 + " Lid"
$ws.Range("AC14").Value = 1

# New row 15: 20 oz All Around (tm) Tumbler bundle test case
$ws.Range("A15").Value = "SKU-T20CPB001"
$ws.Range("AB15").Value = "20 oz All Around" + [char]0x2122 + " Tumbler"
$ws.Range("AC15").Value = 1

# The "Full_RedeemGiftcard" row (old row 36) is now row 38 after the
# insert above; refresh its gift card code.
$ws.Range("O38").Value = "5FD86L34M4337H84S87K"

# Match the author's recorded viewport/selection in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("O27").Select()
